$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value = 7026

$ws1.Range("C3").Value = "合肥·第一届宅舞比赛漫展-CF01（取消）"
$ws1.Range("G3").Value = "不可售"

$ws1.Range("F4").Value = 69
$ws1.Range("F5").Value = 459
$ws1.Range("F6").Value = 158
$ws1.Range("F7").Value = 6958
$ws1.Range("F8").Value = 78
$ws1.Range("F10").Value = 1286
$ws1.Range("F12").Value = 110
$ws1.Range("F13").Value = 412
$ws1.Range("F14").Value = 153
$ws1.Range("F16").Value = 419
$ws1.Range("F18").Value = 48
$ws1.Range("F19").Value = 19
$ws1.Range("F20").Value = 5324
$ws1.Range("F21").Value = 128
$ws1.Range("F22").Value = 183
$ws1.Range("F23").Value = 742
$ws1.Range("F25").Value = 260

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value = 7026

$ws4.Range("C3").Value = "合肥·第一届宅舞比赛漫展-CF01（取消）"
$ws4.Range("G3").Value = "不可售"

$ws4.Range("F4").Value = 69
$ws4.Range("F5").Value = 459
$ws4.Range("F6").Value = 158
$ws4.Range("F7").Value = 6958
$ws4.Range("F8").Value = 78
$ws4.Range("F10").Value = 1286
$ws4.Range("F12").Value = 110
$ws4.Range("F13").Value = 412
$ws4.Range("F14").Value = 153
$ws4.Range("F16").Value = 419
$ws4.Range("F18").Value = 48
$ws4.Range("F19").Value = 19
$ws4.Range("F21").Value = 5324
$ws4.Range("F23").Value = 128
$ws4.Range("F24").Value = 183
$ws4.Range("F25").Value = 742
$ws4.Range("F27").Value = 260

Write-Host "done"
